$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 7.92x33 FMJ (row 4) - price update
$ws.Range("C4").Value = 2250

# 7.62x54 7h1 FMJ (row 6) - price and perf update
$ws.Range("C6").Value = 2250
$ws.Range("G6").Value = 0.3

# 7.62x54 AP (row 7) - price update
$ws.Range("C7").Value = 9500

# Update the selected cell/view, matching the new active selection
$ws.Range("D6").Select()

$wb.Save()
